# Refresh the run's "timestamp" column (Z) on the log sheet with the new
# capture times from this execution. The dataset/feature columns (A:Y) are
# unchanged from the previous run; only Z2:Z48 get new ISO-8601 timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamps = @(
    "2025-10-17T07:09:31.760405",
    "2025-10-17T07:09:31.761405",
    "2025-10-17T07:09:31.761405",
    "2025-10-17T07:09:31.761405",
    "2025-10-17T07:09:31.762406",
    "2025-10-17T07:09:31.762406",
    "2025-10-17T07:09:31.762406",
    "2025-10-17T07:09:31.762406",
    "2025-10-17T07:09:31.763406",
    "2025-10-17T07:09:31.763406",
    "2025-10-17T07:09:31.763406",
    "2025-10-17T07:09:31.763406",
    "2025-10-17T07:09:31.764404",
    "2025-10-17T07:09:31.764404",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.830420",
    "2025-10-17T07:09:31.842252",
    "2025-10-17T07:09:31.842252",
    "2025-10-17T07:09:31.842252",
    "2025-10-17T07:09:31.909592",
    "2025-10-17T07:09:31.909592",
    "2025-10-17T07:09:31.911398",
    "2025-10-17T07:09:31.911398",
    "2025-10-17T07:09:31.911938",
    "2025-10-17T07:09:31.911938",
    "2025-10-17T07:09:31.911938",
    "2025-10-17T07:09:31.911938",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451",
    "2025-10-17T07:09:31.912451"
)

$firstRow = 2
$col = 26  # column Z

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, $col).Value = $timestamps[$i]
}
